# Update the date header and recompute all division answers in the table
# (output regenerated for a new day, per commit "Update master to output generated at 9a8706d")
$d = $word.ActiveDocument

$wdReplaceAll = 2
$wdFindContinue = 1

$found = $d.Content.Find.Execute("2024-02-29 Thursday", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "2024-03-01 Friday", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '2024-02-29 Thursday'" }
$found = $d.Content.Find.Execute("899÷3=299, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "820÷4=205, 0", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '899÷3=299, 2'" }
$found = $d.Content.Find.Execute("162÷7=23, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "740÷9=82, 2", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '162÷7=23, 1'" }
$found = $d.Content.Find.Execute("198÷4=49, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "176÷8=22, 0", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '198÷4=49, 2'" }
$found = $d.Content.Find.Execute("432÷3=144, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "412÷7=58, 6", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '432÷3=144, 0'" }
$found = $d.Content.Find.Execute("324÷5=64, 4", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "644÷9=71, 5", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '324÷5=64, 4'" }
$found = $d.Content.Find.Execute("842÷9=93, 5", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "863÷8=107, 7", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '842÷9=93, 5'" }
$found = $d.Content.Find.Execute("194÷9=21, 5", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "795÷7=113, 4", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '194÷9=21, 5'" }
$found = $d.Content.Find.Execute("308÷2=154, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "601÷6=100, 1", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '308÷2=154, 0'" }
$found = $d.Content.Find.Execute("993÷3=331, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "546÷3=182, 0", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '993÷3=331, 0'" }
$found = $d.Content.Find.Execute("988÷4=247, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "264÷8=33, 0", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '988÷4=247, 0'" }
$found = $d.Content.Find.Execute("896÷7=128, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "534÷8=66, 6", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '896÷7=128, 0'" }
$found = $d.Content.Find.Execute("964÷9=107, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "291÷3=97, 0", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '964÷9=107, 1'" }
$found = $d.Content.Find.Execute("337÷4=84, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "935÷4=233, 3", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '337÷4=84, 1'" }
$found = $d.Content.Find.Execute("340÷9=37, 7", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "544÷6=90, 4", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '340÷9=37, 7'" }
$found = $d.Content.Find.Execute("787÷7=112, 3", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "278÷6=46, 2", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '787÷7=112, 3'" }
$found = $d.Content.Find.Execute("602÷6=100, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "650÷9=72, 2", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '602÷6=100, 2'" }
$found = $d.Content.Find.Execute("998÷2=499, 0", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "345÷9=38, 3", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '998÷2=499, 0'" }
$found = $d.Content.Find.Execute("128÷3=42, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "784÷8=98, 0", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '128÷3=42, 2'" }
$found = $d.Content.Find.Execute("468÷5=93, 3", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "513÷7=73, 2", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '468÷5=93, 3'" }
$found = $d.Content.Find.Execute("867÷7=123, 6", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "499÷3=166, 1", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '867÷7=123, 6'" }
$found = $d.Content.Find.Execute("202÷8=25, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "951÷2=475, 1", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '202÷8=25, 2'" }
$found = $d.Content.Find.Execute("728÷6=121, 2", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "133÷8=16, 5", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '728÷6=121, 2'" }
$found = $d.Content.Find.Execute("233÷6=38, 5", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "683÷8=85, 3", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '233÷6=38, 5'" }
$found = $d.Content.Find.Execute("661÷3=220, 1", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "267÷9=29, 6", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '661÷3=220, 1'" }
$found = $d.Content.Find.Execute("867÷9=96, 3", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "623÷2=311, 1", $wdReplaceAll)
if (-not $found) { Write-Host "WARNING: could not find '867÷9=96, 3'" }
